$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "batsman" column (D), shifting
# batsman/totalRuns/totalBalls/total4s/total6s/sr from D:I to F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data cells for row 2 (existing data row).
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Kings XI Punjab"

# Append a new data row (row 3). Force text storage (matching the rest of
# the sheet, where numeric-looking values are stored as text) by setting
# the number format to Text before assigning the values.
$row3 = $ws.Range("A3:K3")
$row3.NumberFormat = "@"
$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 12 2020"
$ws.Range("C3").Value = "RCB won by 82 runs"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Prasidh Krishna$([char]0x00A0)"
$ws.Range("G3").Value = "2"
$ws.Range("H3").Value = "3"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "66.66"
